# simplifiedPreventiveServices.xlsx — screening logic performance and fixing/testing visits
#
# 1) Six cells in column C ("Frequency") currently hold the text "Once"
#    (rows 2, 3, 20, 21, 23, 26). Replace that value everywhere with "q-1".
#    Three of those cells (rows 20, 21, 23) also carry a slightly different
#    font-color style than the rest of the column; bring them in line with
#    the plain style used by the other "q*" frequency cells (copy formats
#    only from a cell that already has the target style, e.g. C2).
#
# 2) K8 (proportionOfPopulationAtRisk for the first "Skin Cancer Prevention"
#    row) had a stray 0.7 value — clear it back out, keeping its number
#    format.
#
# 3) J30 / J31 ("Healthful Diet and Physical Activity..." Time column) were
#    storing the literal text "12.2 min" — replace with the plain numeric
#    value 12.2 (minutes), matching row 46's "Abnormal Blood Glucose" row
#    which already stores its Time as a bare number.
#
# 4) The sheet's last saved selection moves from L7 to C55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Once" -> "q-1", and normalize style on the three off-style cells ---

$onceStyleFixCells = @("C20", "C21", "C23")
foreach ($addr in $onceStyleFixCells) {
    $ws.Range("C2").Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}

$onceCells = @("C2", "C3", "C20", "C21", "C23", "C26")
foreach ($addr in $onceCells) {
    $ws.Range($addr).Value = "q-1"
}

# --- 2) Clear the stray K8 value ---

$ws.Range("K8").ClearContents()

# --- 3) Replace the "12.2 min" text with the numeric value 12.2 ---

$ws.Range("J30").Value = 12.2
$ws.Range("J31").Value = 12.2

# --- 4) Update the saved selection to C55 ---

$ws.Range("C55").Select()
